$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed Nppc-Npr3 ligand-receptor table (new TPM-derived NATMI run).
# The sending/target cluster set now also includes "ECs" alongside the
# original "FAPs"/"MuSCs", producing a full 3x3 cross of clusters and all
# of the dependent statistics are recomputed accordingly (rows 2-10,
# replacing the previous 6-row table).
$data = @(
    ,("ECs",   "Nppc", "Npr3", "ECs",   1, 0.3333333333333333, 0.01599966666666667, 0.047999,          0.0009948988194802843, 0.0009948988194802843, 3, 1, 0.912887,           2.738661, 0.3341311666818157, 0.3341311666818156, 0.01460588770433333, 0.131452989339,    0.0003324267032833085, 0.0003324267032833085)
    ,("ECs",   "Nppc", "Npr3", "FAPs",  1, 0.3333333333333333, 0.01599966666666667, 0.047999,          0.0009948988194802843, 0.0009948988194802843, 3, 1, 1.302136,           3.906408, 0.4766024939104103, 0.4766024939104103, 0.02083374195466666, 0.187503677592,    0.0004741712585528266, 0.0004741712585528266)
    ,("ECs",   "Nppc", "Npr3", "MuSCs", 1, 0.3333333333333333, 0.01599966666666667, 0.047999,          0.0009948988194802843, 0.0009948988194802843, 3, 1, 0.5170986666666667, 1.551296, 0.1892663394077741, 0.189266339407774,  0.008273406300444444, 0.074460656704,    0.0001883008576441492, 0.0001883008576441492)
    ,("FAPs",  "Nppc", "Npr3", "ECs",   2, 0.6666666666666666, 0.6973903333333333, 2.092171,           0.04336545465636547,   0.04336545465636547,   3, 1, 0.912887,           2.738661, 0.3341311666818157, 0.3341311666818156, 0.6366385692256666,  5.729747123031,    0.01448974995801877,   0.01448974995801877)
    ,("FAPs",  "Nppc", "Npr3", "FAPs",  2, 0.6666666666666666, 0.6973903333333333, 2.092171,           0.04336545465636547,   0.04336545465636547,   3, 1, 1.302136,           3.906408, 0.4766024939104103, 0.4766024939104103, 0.9080970590853333,  8.172873531767999, 0.0206680838387826,    0.0206680838387826)
    ,("FAPs",  "Nppc", "Npr3", "MuSCs", 2, 0.6666666666666666, 0.6973903333333333, 2.092171,           0.04336545465636547,   0.04336545465636547,   3, 1, 0.5170986666666667, 1.551296, 0.1892663394077741, 0.189266339407774,  0.3606196115128889,  3.245576503616,    0.008207620859564103,  0.008207620859564101)
    ,("MuSCs", "Nppc", "Npr3", "ECs",   3, 1,                   15.36831233333333,  46.104937,         0.9556396465241543,    0.9556396465241543,    3, 1, 0.912887,           2.738661, 0.3341311666818157, 0.3341311666818156, 14.02953254103967,   126.265792869357,  0.3193089900205137,    0.3193089900205136)
    ,("MuSCs", "Nppc", "Npr3", "FAPs",  3, 1,                   15.36831233333333,  46.104937,         0.9556396465241543,    0.9556396465241543,    3, 1, 1.302136,           3.906408, 0.4766024939104103, 0.4766024939104103, 20.01163274847733,   180.104694736296,  0.4554602388130749,    0.4554602388130749)
    ,("MuSCs", "Nppc", "Npr3", "MuSCs", 3, 1,                   15.36831233333333,  46.104937,         0.9556396465241543,    0.9556396465241543,    3, 1, 0.5170986666666667, 1.551296, 0.1892663394077741, 0.189266339407774,  7.946933816483557,   71.522404348352,   0.1808704176905658,    0.1808704176905658)
)

$rowIndex = 2
foreach ($row in $data) {
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}
